$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the formatting of the
# neighboring header cell (G1) so it matches the other header cells exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2 (plain, unstyled numeric cell)
$ws.Range("H2").Value = 0
